$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> Adam10 -> Epha3 -> ECs
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 49.46095133333333
$ws.Range("H2").Value = 148.382854
$ws.Range("I2").Value = 0.5804304915153436
$ws.Range("J2").Value = 0.5804304915153436
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.003058333333333333
$ws.Range("N2").Value = 0.009175000000000001
$ws.Range("O2").Value = 0.0001134234803787887
$ws.Range("P2").Value = 0.0001134234803787887
$ws.Range("Q2").Value = 0.1512680761611111
$ws.Range("R2").Value = 1.36141268545
$ws.Range("S2").Value = 0.00006583444646564125
$ws.Range("T2").Value = 0.00006583444646564125

# Row 3: ECs -> Adam10 -> Epha3 -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 49.46095133333333
$ws.Range("H3").Value = 148.382854
$ws.Range("I3").Value = 0.5804304915153436
$ws.Range("J3").Value = 0.5804304915153436
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 26.097779
$ws.Range("N3").Value = 78.29333700000001
$ws.Range("O3").Value = 0.9678804112271815
$ws.Range("P3").Value = 0.9678804112271815
$ws.Range("Q3").Value = 1290.820977027089
$ws.Range("R3").Value = 11617.3887932438
$ws.Range("S3").Value = 0.5617873028166659
$ws.Range("T3").Value = 0.5617873028166659

# Row 4: ECs -> Adam10 -> Epha3 -> sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 49.46095133333333
$ws.Range("H4").Value = 148.382854
$ws.Range("I4").Value = 0.5804304915153436
$ws.Range("J4").Value = 0.5804304915153436
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.8630093333333333
$ws.Range("N4").Value = 2.589028
$ws.Range("O4").Value = 0.03200616529243972
$ws.Range("P4").Value = 0.03200616529243972
$ws.Range("Q4").Value = 42.68526263621244
$ws.Range("R4").Value = 384.167363725912
$ws.Range("S4").Value = 0.01857735425221212
$ws.Range("T4").Value = 0.01857735425221212

# Row 5: FAPs -> Adam10 -> Epha3 -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Adam10"
$ws.Range("C5").Value = "Epha3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 18.03569333333333
$ws.Range("H5").Value = 54.10708
$ws.Range("I5").Value = 0.2116511321372752
$ws.Range("J5").Value = 0.2116511321372752
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.003058333333333333
$ws.Range("N5").Value = 0.009175000000000001
$ws.Range("O5").Value = 0.0001134234803787887
$ws.Range("P5").Value = 0.0001134234803787887
$ws.Range("Q5").Value = 0.0551591621111111
$ws.Range("R5").Value = 0.496432459
$ws.Range("S5").Value = 0.00002400620803312064
$ws.Range("T5").Value = 0.00002400620803312064

# Row 6: FAPs -> Adam10 -> Epha3 -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Adam10"
$ws.Range("C6").Value = "Epha3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 18.03569333333333
$ws.Range("H6").Value = 54.10708
$ws.Range("I6").Value = 0.2116511321372752
$ws.Range("J6").Value = 0.2116511321372752
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 26.097779
$ws.Range("N6").Value = 78.29333700000001
$ws.Range("O6").Value = 0.9678804112271815
$ws.Range("P6").Value = 0.9678804112271815
$ws.Range("Q6").Value = 470.6915387251066
$ws.Range("R6").Value = 4236.22384852596
$ws.Range("S6").Value = 0.2048529848097245
$ws.Range("T6").Value = 0.2048529848097245

# Row 7: FAPs -> Adam10 -> Epha3 -> sCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Adam10"
$ws.Range("C7").Value = "Epha3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 18.03569333333333
$ws.Range("H7").Value = 54.10708
$ws.Range("I7").Value = 0.2116511321372752
$ws.Range("J7").Value = 0.2116511321372752
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.8630093333333333
$ws.Range("N7").Value = 2.589028
$ws.Range("O7").Value = 0.03200616529243972
$ws.Range("P7").Value = 0.03200616529243972
$ws.Range("Q7").Value = 15.56497167980444
$ws.Range("R7").Value = 140.08474511824
$ws.Range("S7").Value = 0.006774141119517631
$ws.Range("T7").Value = 0.006774141119517631

# Row 8: sCs -> Adam10 -> Epha3 -> ECs
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Adam10"
$ws.Range("C8").Value = "Epha3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 17.71760933333333
$ws.Range("H8").Value = 53.152828
$ws.Range("I8").Value = 0.2079183763473812
$ws.Range("J8").Value = 0.2079183763473812
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.003058333333333333
$ws.Range("N8").Value = 0.009175000000000001
$ws.Range("O8").Value = 0.0001134234803787887
$ws.Range("P8").Value = 0.0001134234803787887
$ws.Range("Q8").Value = 0.05418635521111111
$ws.Range("R8").Value = 0.4876771969000001
$ws.Range("S8").Value = 0.00002358282588002679
$ws.Range("T8").Value = 0.00002358282588002679

# Row 9: sCs -> Adam10 -> Epha3 -> FAPs
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Adam10"
$ws.Range("C9").Value = "Epha3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 17.71760933333333
$ws.Range("H9").Value = 53.152828
$ws.Range("I9").Value = 0.2079183763473812
$ws.Range("J9").Value = 0.2079183763473812
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 26.097779
$ws.Range("N9").Value = 78.29333700000001
$ws.Range("O9").Value = 0.9678804112271815
$ws.Range("P9").Value = 0.9678804112271815
$ws.Range("Q9").Value = 462.3902527896707
$ws.Range("R9").Value = 4161.512275107037
$ws.Range("S9").Value = 0.2012401236007912
$ws.Range("T9").Value = 0.2012401236007912

# Row 10: sCs -> Adam10 -> Epha3 -> sCs
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Adam10"
$ws.Range("C10").Value = "Epha3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 17.71760933333333
$ws.Range("H10").Value = 53.152828
$ws.Range("I10").Value = 0.2079183763473812
$ws.Range("J10").Value = 0.2079183763473812
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.8630093333333333
$ws.Range("N10").Value = 2.589028
$ws.Range("O10").Value = 0.03200616529243972
$ws.Range("P10").Value = 0.03200616529243972
$ws.Range("Q10").Value = 15.29046221902044
$ws.Range("R10").Value = 137.614159971184
$ws.Range("S10").Value = 0.006654669920709972
$ws.Range("T10").Value = 0.006654669920709972
